# Generate Report for Handback
# Update the "latest generated" timestamps for the 7910e7a4-... file entry
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 corresponds to the 7910e7a4-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-08 05:05:17"

# --- zh-cn sheet: row 3 corresponds to the 7910e7a4-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-08 05:05:12"
$wsZhCn.Range("K3").Value = "2016-09-08 05:05:29"

# --- de-de sheet: row 3 corresponds to the 7910e7a4-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-08 05:05:17"
$wsDeDe.Range("K3").Value = "2016-09-08 05:05:38"
